# Update Active_Outages.xlsx - 6/18/2025, 5:02:07 PM
#
# Refreshes the "Elapsed Duration(Hrs)" counters on every region sheet and
# updates/adds a couple of outage rows (new JED0125 -> JED0190 relocation
# on R1/R2, and a new MAK0875 entry on R5).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet R1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("R1")

# Refresh elapsed-duration counters for the existing outages.
$ws1.Range("G2").Value = "3930:16:06"
$ws1.Range("G3").Value = "69:48:44"
$ws1.Range("G4").Value = "92:48:44"

# Row 5 site moved from region R5 / HAJ0155 to region R4 / JED0125 and the
# power source / battery-backup feedback were updated accordingly.
$ws1.Range("B5").Value = "R4"
$ws1.Range("D5").Value = "JED0125"
$ws1.Range("I5").Value = "Generator-SG"
$ws1.Range("J5").Value = "Good+In progress"

# ---------------------------------------------------------------------
# Sheet R2
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("R2")

$ws2.Range("G2").Value = "12111:39:47"
$ws2.Range("G3").Value = "3241:23:16"
$ws2.Range("G4").Value = "479:34:50"

# New outage row added to R2.
$ws2.Range("B5").Value = "R4"
$ws2.Range("D5").Value = "JED0190"
$ws2.Range("I5").Value = "SCECO"
$ws2.Range("J5").Value = "Good"
$ws2.Range("L5").Value = "Latis"

# ---------------------------------------------------------------------
# Sheet R4
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("R4")

$ws4.Range("G2").Value = "2957:29:36"
$ws4.Range("G3").Value = "184:41:51"
$ws4.Range("G4").Value = "72:54:16"
$ws4.Range("G5").Value = "70:31:49"

# ---------------------------------------------------------------------
# Sheet R5
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("R5")

$ws5.Range("G2").Value = "431:28:35"

# New outage row added to R5 (site only, rest of the row left blank).
$ws5.Range("A3").Value = "MAK0875"

# ---------------------------------------------------------------------
# Sheet R6
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("R6")

$ws6.Range("G2").Value = "72:00:53"
